# Re-ran "resolve" and "classify + summarise" steps after changes to the
# mapping file. For Nagaland, this reduced the set of species feeding the
# Range Status / Species qualification / High Priority break-up summaries.

$wb = $excel.ActiveWorkbook

# --- Sheet "Range Status": zero-out species counts (col B), drop percentages (col C) ---
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("B3").Value = 0
$wsRange.Range("B4").Value = 0
$wsRange.Range("B5").Value = 0
$wsRange.Range("B6").Value = 0
$wsRange.Range("B7").Value = 0
$wsRange.Range("C2:C7").ClearContents()

# --- Sheet "Species qualification": Range Analysis row now has 0 species ---
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# --- Sheet "High Priority break-up": only the IUCN break-up remains, and
#     it now accounts for all 7 high-priority species ---
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Rows("3").Delete()
$wsBreak.Range("A2").Value = "IUCN"
$wsBreak.Range("B2").Value = 7
$wsBreak.Range("C2").Value = 100
$wsBreak.Range("D2").Value = 7
$wsBreak.Range("E2").Value = 100
